$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (N) that mirrors the existing "2019" column (M):
# same formatting as M4/M5, new values 2020 / 2.1.
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N4").Value = 2020

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N5").Value = 2.1

$excel.CutCopyMode = $false

# Match the saved selection state (active cell N9) recorded in the file.
$ws.Range("N9").Select()
